# Update column G ("K") values for rows 2-29 on the active sheet.
# These are the new "K" values calculated to replace the old "Strike#" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 4
    3  = 6
    4  = 7
    5  = 8
    6  = 2
    7  = 1
    8  = 4
    9  = 8
    10 = 2
    11 = 2
    12 = 10
    13 = 3
    14 = 5
    15 = 6
    16 = 3
    17 = 1
    18 = 4
    19 = 2
    20 = 8
    21 = 7
    22 = 3
    23 = 6
    24 = 3
    25 = 4
    26 = 4
    27 = 5
    28 = 2
    29 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
